$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Datos crudos")
$ws2 = $wb.Worksheets.Item("Datos válidos")

$newData = @(
    @("2023-12-11 23:54:44", 25.935314685314701),
    @("2023-12-11 23:55:45", 25.1486013986014),
    @("2023-12-11 23:56:47", 24.230769230769202),
    @("2023-12-11 23:57:49", 22.788461538461501),
    @("2023-12-11 23:58:51", 22.6573426573426),
    @("2023-12-11 23:59:53", 22.263986013985999),
    @("2023-12-12 00:00:55", 22.0017482517482),
    @("2023-12-12 00:01:57", 22.0017482517482),
    @("2023-12-12 00:02:59", 21.608391608391599),
    @("2023-12-12 00:04:01", 21.346153846153801),
    @("2023-12-12 00:05:03", 21.870629370629299),
    @("2023-12-12 00:06:05", 21.346153846153801),
    @("2023-12-12 00:07:07", 21.215034965034899),
    @("2023-12-12 00:08:09", 21.346153846153801),
    @("2023-12-12 00:09:11", 21.083916083916101),
    @("2023-12-12 00:10:13", 20.821678321678299),
    @("2023-12-12 00:11:15", 21.083916083916101),
    @("2023-12-12 00:12:17", 21.215034965034899),
    @("2023-12-12 00:13:19", 21.215034965034899),
    @("2023-12-12 00:14:21", 21.346153846153801),
    @("2023-12-12 00:15:23", 21.083916083916101),
    @("2023-12-12 00:16:25", 20.821678321678299),
    @("2023-12-12 00:17:27", 20.821678321678299),
    @("2023-12-12 00:18:29", 21.477272727272702),
    @("2023-12-12 00:19:31", 21.346153846153801),
    @("2023-12-12 00:20:33", 20.821678321678299),
    @("2023-12-12 00:21:35", 20.5594405594405),
    @("2023-12-12 00:22:37", 21.346153846153801),
    @("2023-12-12 00:23:39", 21.083916083916101),
    @("2023-12-12 00:24:41", 21.083916083916101),
    @("2023-12-12 00:25:43", 20.690559440559401),
    @("2023-12-12 00:26:45", 20.821678321678299),
    @("2023-12-12 00:27:47", 21.083916083916101),
    @("2023-12-12 00:28:49", 20.821678321678299),
    @("2023-12-12 00:29:51", 21.083916083916101),
    @("2023-12-12 00:30:53", 20.9527972027972),
    @("2023-12-12 00:31:55", 20.821678321678299),
    @("2023-12-12 00:32:57", 21.346153846153801),
    @("2023-12-12 00:33:58", 22.0017482517482)
)

$dateFormat = $ws1.Range("C37").NumberFormat

for ($i = 0; $i -lt $newData.Count; $i++) {
    $row = 2 + $i
    $ts = $newData[$i][0]
    $temp = $newData[$i][1]

    $ws1.Cells.Item($row, 1).Value = 23
    $ws1.Cells.Item($row, 2).Value = 19
    $c = $ws1.Cells.Item($row, 3)
    $c.NumberFormat = $dateFormat
    $c.Value = $ts
    $ws1.Cells.Item($row, 4).Value = 0
    $ws1.Cells.Item($row, 5).Value = $temp
}

$ws1.Range("H3").Formula = "=COUNT(E:E)-2"

$ws2.Range("B3").Select()
$ws1.Activate()
